# Updated symbol list on Thu Feb 16 01:07:15 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D), Volume(1h) (E), Data (F) and Hora (G) columns for
# every coin row (2-51) on the active sheet to the next scrape snapshot.
#
# D/E/G hold numeric-looking text ("320.35", "8.57%", "1") that must stay
# plain text (the sheet already stores Price/Volume/Hora as inlineStr, not
# numbers/dates), so those values are written with a leading apostrophe to
# force Excel's text interpretation instead of auto-converting them to a
# number/percentage. F holds a literal "D-M-YYYY" string which Excel already
# keeps as text as-is, so no apostrophe is needed there. Rows whose Price or
# Volume did not change (still "--"/"--%", i.e. no market data) are left
# untouched via $null.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "320.35"; E = "8.57%"; F = "16-2-2023"; G = "1" },
    @{ Row = 3; D = "45.14"; E = "7.09%"; F = "16-2-2023"; G = "1" },
    @{ Row = 4; D = "5.185"; E = "3.34%"; F = "16-2-2023"; G = "1" },
    @{ Row = 5; D = "0.08081"; E = "7.31%"; F = "16-2-2023"; G = "1" },
    @{ Row = 6; D = "4.579"; E = "4.15%"; F = "16-2-2023"; G = "1" },
    @{ Row = 7; D = "1.668"; E = "5.10%"; F = "16-2-2023"; G = "1" },
    @{ Row = 8; D = "1.091"; E = "17.57%"; F = "16-2-2023"; G = "1" },
    @{ Row = 9; D = "0.1327"; E = "10.42%"; F = "16-2-2023"; G = "1" },
    @{ Row = 10; D = "0.1947"; E = "5.85%"; F = "16-2-2023"; G = "1" },
    @{ Row = 11; D = "0.09437"; E = "5.78%"; F = "16-2-2023"; G = "1" },
    @{ Row = 12; D = "0.04351"; E = "6.89%"; F = "16-2-2023"; G = "1" },
    @{ Row = 13; D = "0.1040"; E = "-1.01%"; F = "16-2-2023"; G = "1" },
    @{ Row = 14; D = $null; E = "2.31%"; F = "16-2-2023"; G = "1" },
    @{ Row = 15; D = "0.005785"; E = "-0.43%"; F = "16-2-2023"; G = "1" },
    @{ Row = 16; D = $null; E = $null; F = "16-2-2023"; G = "1" },
    @{ Row = 17; D = "3.423"; E = "1.83%"; F = "16-2-2023"; G = "1" },
    @{ Row = 18; D = "2.426"; E = "0.27%"; F = "16-2-2023"; G = "1" },
    @{ Row = 19; D = "0.3383"; E = "1.89%"; F = "16-2-2023"; G = "1" },
    @{ Row = 20; D = "8.261"; E = "5.37%"; F = "16-2-2023"; G = "1" },
    @{ Row = 21; D = "0.1357"; E = "-1.84%"; F = "16-2-2023"; G = "1" },
    @{ Row = 22; D = $null; E = "4.76%"; F = "16-2-2023"; G = "1" },
    @{ Row = 23; D = "0.04288"; E = "5.44%"; F = "16-2-2023"; G = "1" },
    @{ Row = 24; D = "0.001295"; E = "2.26%"; F = "16-2-2023"; G = "1" },
    @{ Row = 25; D = "0.004255"; E = "5.27%"; F = "16-2-2023"; G = "1" },
    @{ Row = 26; D = "0.0001344"; E = "9.12%"; F = "16-2-2023"; G = "1" },
    @{ Row = 27; D = $null; E = $null; F = "16-2-2023"; G = "1" },
    @{ Row = 28; D = $null; E = $null; F = "16-2-2023"; G = "1" },
    @{ Row = 29; D = $null; E = $null; F = "16-2-2023"; G = "1" },
    @{ Row = 30; D = $null; E = $null; F = "16-2-2023"; G = "1" },
    @{ Row = 31; D = $null; E = $null; F = "16-2-2023"; G = "1" },
    @{ Row = 32; D = $null; E = $null; F = "16-2-2023"; G = "1" },
    @{ Row = 33; D = $null; E = $null; F = "16-2-2023"; G = "1" },
    @{ Row = 34; D = $null; E = $null; F = "16-2-2023"; G = "1" },
    @{ Row = 35; D = $null; E = $null; F = "16-2-2023"; G = "1" },
    @{ Row = 36; D = $null; E = $null; F = "16-2-2023"; G = "1" },
    @{ Row = 37; D = $null; E = $null; F = "16-2-2023"; G = "1" },
    @{ Row = 38; D = "0.02676"; E = "11.41%"; F = "16-2-2023"; G = "1" },
    @{ Row = 39; D = "0.05468"; E = "5.36%"; F = "16-2-2023"; G = "1" },
    @{ Row = 40; D = "0.005860"; E = "-3.29%"; F = "16-2-2023"; G = "1" },
    @{ Row = 41; D = "0.007773"; E = "-0.27%"; F = "16-2-2023"; G = "1" },
    @{ Row = 42; D = "0.1438"; E = "8.28%"; F = "16-2-2023"; G = "1" },
    @{ Row = 43; D = "0.007356"; E = "-2.45%"; F = "16-2-2023"; G = "1" },
    @{ Row = 44; D = "0.008565"; E = "18.18%"; F = "16-2-2023"; G = "1" },
    @{ Row = 45; D = "0.3188"; E = "-0.38%"; F = "16-2-2023"; G = "1" },
    @{ Row = 46; D = "0.00006854"; E = "1.50%"; F = "16-2-2023"; G = "1" },
    @{ Row = 47; D = $null; E = "-0.43%"; F = "16-2-2023"; G = "1" },
    @{ Row = 48; D = "0.05340"; E = "32.06%"; F = "16-2-2023"; G = "1" },
    @{ Row = 49; D = "0.003988"; E = "-5.21%"; F = "16-2-2023"; G = "1" },
    @{ Row = 50; D = "0.00002094"; E = "-0.43%"; F = "16-2-2023"; G = "1" },
    @{ Row = 51; D = $null; E = "-0.43%"; F = "16-2-2023"; G = "1" }
)

foreach ($u in $updates) {
    # Force-text (leading apostrophe): Price, Volume(1h) and Hora look numeric.
    if ($null -ne $u.D) { $ws.Range("D" + $u.Row).Value = "'" + $u.D }
    if ($null -ne $u.E) { $ws.Range("E" + $u.Row).Value = "'" + $u.E }
    # Data is a "D-M-YYYY" style string; Excel keeps it as text natively.
    if ($null -ne $u.F) { $ws.Range("F" + $u.Row).Value = $u.F }
    if ($null -ne $u.G) { $ws.Range("G" + $u.Row).Value = "'" + $u.G }
}
